$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column H, matching the formatting of the existing
# header cells (B1:G1) by copying G1's format onto H1, then overwriting
# the value.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data cells H2:H5 = 0, H6 = 1 (plain, unstyled numbers like column F/G)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
